$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the current last sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Match the page margins used throughout the rest of the workbook
$newSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.TopMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# These columns hold numeric-looking text values (match codes, counts, percentages,
# "NO") that must stay literal text rather than being auto-converted to numbers.
$newSheet.Range("A2:A4").NumberFormat = "@"
$newSheet.Range("C2:F4").NumberFormat = "@"

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Data rows
$newSheet.Range("A2").Value = "4452"
$newSheet.Range("B2").Value = 8
$newSheet.Range("C2").Value = "1"
$newSheet.Range("D2").Value = "1"
$newSheet.Range("E2").Value = "10.69%"
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4453"
$newSheet.Range("B3").Value = 7
$newSheet.Range("C3").Value = "0"
$newSheet.Range("D3").Value = "1"
$newSheet.Range("E3").Value = "2.58%"
$newSheet.Range("F3").Value = "NO"

$newSheet.Range("A4").Value = "4455"
$newSheet.Range("B4").Value = 8
$newSheet.Range("C4").Value = "0"
$newSheet.Range("D4").Value = "0"
$newSheet.Range("E4").Value = "1.95%"
$newSheet.Range("F4").Value = "NO"

# Give the header row the same bold / bordered / centered look used by the other
# sheets' header rows, by copying the formatting from an existing header cell so
# the same underlying style gets reused.
$srcSheet = $wb.Worksheets.Item("ODI Bowling")
$srcSheet.Range("A1").Copy() | Out-Null
$newSheet.Range("A1:F1").PasteSpecial(-4122) | Out-Null
